# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the 78d0cb2c-0ee6-4e4b-a55b-b4f49e5266fb.md row (row 7
# in the table / row 6 of the sheetData, i.e. the 5th data row) to reflect
# a freshly generated handoff package.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G6").Value = "2016-08-24 16:43:23"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H6").Value = "2016-08-24 16:43:18"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H6").Value = "2016-08-24 16:43:23"
